$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All "Price" column (D) values are stored as text in this sheet (the source
# file uses inline strings everywhere). Force text formatting before writing
# numeric-looking strings so Excel doesn't silently convert them to numbers.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Simple price updates (column D) for rows where only the price changed
Set-TextValue "D2" "249.27"
Set-TextValue "D3" "22.03"
Set-TextValue "D4" "5.346"
Set-TextValue "D5" "0.05628"
Set-TextValue "D7" "6.377"
Set-TextValue "D8" "0.8169"
Set-TextValue "D9" "0.9279"
Set-TextValue "D10" "0.1448"
Set-TextValue "D11" "0.07448"
Set-TextValue "D12" "0.03242"
Set-TextValue "D13" "0.03083"
Set-TextValue "D14" "0.09318"
Set-TextValue "D15" "3.555"
Set-TextValue "D16" "0.001596"
Set-TextValue "D17" "0.04741"

# Rows 18-24: coin list shifted by one position (new coin "One" inserted at
# row 18, pushing the others down one row, with the previous row-24 coin
# "BTSEToken" now landing at row 24 too, with an updated price/rank label)
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D18" "0.0005758"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D19" "0.006384"
$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D20" "0.005068"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D21" "0.001034"
$ws.Range("E21").Value = "20BitKanKAN"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D22" "0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D23" "3.737"
$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D24" "2.161"
$ws.Range("E24").Value = "23BTSETokenBTSE"

# Row 27: remove "Worstin24h" suffix from the volume label
$ws.Range("E27").Value = "26AAXTokenAAB"

# Rows 40-50 price / label updates
Set-TextValue "D40" "0.03941"
Set-TextValue "D41" "0.002917"
$ws.Range("E41").Value = "40KickTokenKICKWorstin24h"
Set-TextValue "D42" "0.1069"
Set-TextValue "D43" "0.002859"
Set-TextValue "D44" "0.008519"
Set-TextValue "D45" "0.00005570"
Set-TextValue "D49" "0.1928"
Set-TextValue "D50" "0.00002099"
